# Applies the numeric corrections from the scheduled-runner update.
# For each touched row, sets the new cell values and clears any cells
# that were removed entirely (e.g. N column dropped when a leve no longer
# has an HQ variant).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1450
$ws.Range("I12").Value = 1450
$ws.Range("K12").Value = 1450
$ws.Range("M12").Value = -1280
$ws.Range("H33").Value = 16667038
$ws.Range("I33").Value = 17857452
$ws.Range("J33").Value = 1250
$ws.Range("K33").Value = 17857452
$ws.Range("L33").Value = 1250
$ws.Range("M33").Value = -17857223
$ws.Range("N33").Value = -1708
$ws.Range("H38").Value = 555702.2
$ws.Range("I38").Value = 555702.2
$ws.Range("K38").Value = 1667106.6
$ws.Range("M38").Value = -1666734.6
$ws.Range("H64").Value = 5120.222
$ws.Range("I64").Value = 4680.8335
$ws.Range("J64").Value = 5999
$ws.Range("K64").Value = 4680.8335
$ws.Range("L64").Value = 5999
$ws.Range("M64").Value = -4432.8335
$ws.Range("N64").Value = -6495
$ws.Range("H67").Value = 5120.222
$ws.Range("I67").Value = 4680.8335
$ws.Range("J67").Value = 5999
$ws.Range("K67").Value = 4680.8335
$ws.Range("L67").Value = 5999
$ws.Range("M67").Value = -3822.8335
$ws.Range("N67").Value = -7715
$ws.Range("H70").Value = 1644.4286
$ws.Range("I70").Value = 1601.8334
$ws.Range("J70").Value = 1900
$ws.Range("K70").Value = 4805.5002
$ws.Range("L70").Value = 5700
$ws.Range("M70").Value = -4535.5002
$ws.Range("N70").Value = -6240
$ws.Range("H73").Value = 1644.4286
$ws.Range("I73").Value = 1601.8334
$ws.Range("J73").Value = 1900
$ws.Range("K73").Value = 4805.5002
$ws.Range("L73").Value = 5700
$ws.Range("M73").Value = -3869.5002
$ws.Range("N73").Value = -7572
$ws.Range("H74").Value = 10209.875
$ws.Range("I74").Value = 10097
$ws.Range("K74").Value = 10097
$ws.Range("M74").Value = -9161
$ws.Range("H77").Value = 10209.875
$ws.Range("I77").Value = 10097
$ws.Range("K77").Value = 50485
$ws.Range("M77").Value = -45805
$ws.Range("H112").Value = 2800.55
$ws.Range("J112").Value = 2800.55
$ws.Range("L112").Value = 8401.650000000001
$ws.Range("N112").Value = -10617.65
$ws.Range("H135").Value = 1315.0834
$ws.Range("I135").Value = 1178.2
$ws.Range("K135").Value = 10603.8
$ws.Range("M135").Value = -8068.800000000001
$ws.Range("H137").Value = 20835724
$ws.Range("I137").Value = 66668210
$ws.Range("K137").Value = 200004630
$ws.Range("M137").Value = -200002080
$ws.Range("H141").Value = 2897
$ws.Range("I141").Value = 2628.5
$ws.Range("J141").Value = 3434
$ws.Range("K141").Value = 7885.5
$ws.Range("L141").Value = 10302
$ws.Range("M141").Value = -2705.5
$ws.Range("N141").Value = -20662

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1691.8572
$ws.Range("I45").Value = 1691.8572
$ws.Range("K45").Value = 1691.8572
$ws.Range("M45").Value = -1314.8572
$ws.Range("H74").Value = 4602968.5
$ws.Range("I74").Value = 5752461
$ws.Range("K74").Value = 5752461
$ws.Range("M74").Value = -5751587
$ws.Range("H77").Value = 4602968.5
$ws.Range("I77").Value = 5752461
$ws.Range("K77").Value = 28762305
$ws.Range("M77").Value = -28757937

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 771.1429000000001
$ws.Range("I22").Value = 719.6
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 719.6
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -546.6
$ws.Range("N22").Value = -1246
$ws.Range("H35").Value = 39999.5
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H86").Value = 10461518
$ws.Range("I86").Value = 57826.188
$ws.Range("K86").Value = 57826.188
$ws.Range("M86").Value = -56703.188
$ws.Range("H89").Value = 10461518
$ws.Range("I89").Value = 57826.188
$ws.Range("K89").Value = 289130.94
$ws.Range("M89").Value = -283514.94
$ws.Range("H134").Value = 1060.3334
$ws.Range("I134").Value = 872.5
$ws.Range("J134").Value = 1999.5
$ws.Range("K134").Value = 2617.5
$ws.Range("L134").Value = 5998.5
$ws.Range("M134").Value = -82.5
$ws.Range("N134").Value = -11068.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 12110.429
$ws.Range("I99").Value = 3212.1667
$ws.Range("K99").Value = 3212.1667
$ws.Range("M99").Value = -1714.1667
$ws.Range("H122").Value = 3923.2727
$ws.Range("I122").Value = 4015.7
$ws.Range("K122").Value = 12047.1
$ws.Range("M122").Value = -9597.099999999999
$ws.Range("H126").Value = 12110.429
$ws.Range("I126").Value = 3212.1667
$ws.Range("K126").Value = 9636.500100000001
$ws.Range("M126").Value = -7166.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 831
$ws.Range("I64").Value = 831
$ws.Range("K64").Value = 2493
$ws.Range("M64").Value = -2223
$ws.Range("H67").Value = 831
$ws.Range("I67").Value = 831
$ws.Range("K67").Value = 2493
$ws.Range("M67").Value = -1557
$ws.Range("H131").Value = 1384
$ws.Range("J131").Value = 1365.6666
$ws.Range("L131").Value = 4096.9998
$ws.Range("N131").Value = -14176.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3024.8333
$ws.Range("I80").Value = 2673.5
$ws.Range("K80").Value = 2673.5
$ws.Range("M80").Value = -1675.5
$ws.Range("H83").Value = 3024.8333
$ws.Range("I83").Value = 2673.5
$ws.Range("K83").Value = 13367.5
$ws.Range("M83").Value = -8375.5
$ws.Range("H126").Value = 3108
$ws.Range("I126").Value = 2233
$ws.Range("J126").Value = 3983
$ws.Range("K126").Value = 6699
$ws.Range("L126").Value = 11949
$ws.Range("M126").Value = -4229
$ws.Range("N126").Value = -16889

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 14899
$ws.Range("I7").Value = 14899
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 14899
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -14787
$ws.Range("N7").ClearContents()
$ws.Range("H40").Value = 35722296
$ws.Range("I40").Value = 50007400
$ws.Range("K40").Value = 50007400
$ws.Range("M40").Value = -50007264
$ws.Range("H82").Value = 1825.2727
$ws.Range("I82").Value = 1625
$ws.Range("J82").Value = 2065.6
$ws.Range("K82").Value = 1625
$ws.Range("L82").Value = 2065.6
$ws.Range("M82").Value = -1264
$ws.Range("N82").Value = -2787.6
$ws.Range("H85").Value = 1825.2727
$ws.Range("I85").Value = 1625
$ws.Range("J85").Value = 2065.6
$ws.Range("K85").Value = 1625
$ws.Range("L85").Value = 2065.6
$ws.Range("M85").Value = -377
$ws.Range("N85").Value = -4561.6
$ws.Range("H98").Value = 64498.5
$ws.Range("J98").Value = 64498.5
$ws.Range("L98").Value = 64498.5
$ws.Range("N98").Value = -70488.5
$ws.Range("H122").Value = 3850
$ws.Range("I122").Value = 3466.6667
$ws.Range("K122").Value = 10400.0001
$ws.Range("M122").Value = -7950.000100000001
$ws.Range("H126").Value = 14899
$ws.Range("I126").Value = 14899
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 44697
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -42227
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 16666
$ws.Range("I132").Value = 5554.6665
$ws.Range("K132").Value = 16663.9995
$ws.Range("M132").Value = -14133.9995
$ws.Range("H136").Value = 2507.2666
$ws.Range("I136").Value = 2479.2144
$ws.Range("K136").Value = 7437.6432
$ws.Range("M136").Value = -4887.6432
